# Algebra3_Day_006 1.2 Real Number Intro.pptx
# Commit: "edit powerpoints, add HLQ - 2018"
#
# 1) The cached "datetimeFigureOut" footer field (used by the slide master
#    and every slide layout) gets refreshed from 9/5/2017 -> 9/17/2018.
# 2) On slide 1, the two adjacent runs "Graph " and "the following: "
#    (after the title's line break) are merged into a single run
#    "Graph the following: ".

$p = $ppt.ActivePresentation

$oldDate = "9/5/2017"
$newDate = "9/17/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isPlaceholder = $false
        $phType = -1
        try {
            $phType = $sh.PlaceholderFormat.Type
            $isPlaceholder = $true
        } catch {
            $isPlaceholder = $false
        }
        if ($isPlaceholder -and $phType -eq 16) {
            if ($sh.HasTextFrame) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own Date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: merge the "Graph " / "the following: " runs into one run,
# leaving the preceding "Bell Work" line (and its line break) untouched.
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$fullText = $tr.Text
$marker = "Graph the following: "
$startIdx = $fullText.IndexOf($marker)
if ($startIdx -ge 0) {
    $rangeLen = $marker.Length
    # Re-assigning the exact same text is a no-op in this engine, so first
    # nudge the range to a different value, then set the final merged text;
    # this collapses the two runs into a single run (the surrounding
    # "Bell Work" run and the <a:br> line break are left alone).
    $chars = $tr.Characters($startIdx + 1, $rangeLen)
    $chars.Text = $marker + [char]0x2060
    $chars2 = $tr.Characters($startIdx + 1, $rangeLen + 1)
    $chars2.Text = $marker
}
